$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Paragraph 10: "ricev.: Ingegneria Sede Scientifica – Palazzina 1 – giovedì 11:00-12:00"
# Update the second run's text (office hours day/time).
$paraRicev = $tr.Paragraphs(10, 1)
$runRicev = $paraRicev.Runs(2, 1)
$runRicev.Text = ".: Ingegneria Sede Scientifica – Palazzina 1 – lunedì 15:00-16:00"

# Paragraph 11: "materiale: http://elly.dia.unipr.it/2017/course/view.php?id=417"
# Update the third run's text (course material link).
$paraMateriale = $tr.Paragraphs(11, 1)
$runLink = $paraMateriale.Runs(3, 1)
$runLink.Text = "https://albertoferrari.github.io/"

# Remove the now-superfluous trailing empty paragraph (paragraph 12).
$paraTrailing = $tr.Paragraphs(12, 1)
$paraTrailing.Delete()
